$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.315, 0.279, -0.479, -61, -65, -112),
    @(0.011, 0.141, 0.084, -70, 8, -71),
    @(0.098, -0.171, 0.381, -67, -16, -139),
    @(0.521, 0.294, -0.073, 80, -43, -156),
    @(0.073, -0.171, 0.869, -122, -15, -41),
    @(0.62, 0.08, 0.57, -57, 25, 22),
    @(0.067, 0.223, 0.019, 129, -21, 80),
    @(0.167, -0.084, 0.587, -155, 50, 174),
    @(0.62, 0.08, 0.57, -57, 25, 22),
    @(-0.197, 0.113, -0.351, -97, -44, 43)
)

for ($r = 1; $r -le 10; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$r - 1][$c - 1]
    }
}
